$d = $word.ActiveDocument

# 1) Move the "_GoBack" bookmark from its old location (end of the
#    "...reg@fjk.hu ... cimre." paragraph) to a new spot in the middle of
#    the "...szolgaltatasok..." run, splitting that run in two.
$r = $d.Content
$r.Find.Execute("szolgáltatások, a szállást", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $r.Start + ("szolgáltatáso").Length
$insertRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $insertRange) | Out-Null

# 2) Collapse the long bold red paragraph about foreign bank transfers down
#    to a single simplified sentence/run.
$oldText = "Amennyiben az utalást csak külföldi (nem magyar) bankszámláról tudnád intézni, úgy nem vagy köteles az előleg fizetésére tekintettel a magas határon kívüli utalási költségekre. A szállás költségét a Benczúr Hotelben tudod majd egyenlíteni, a regisztráció fennmaradó költséget pedig a konferencia regisztrációs asztalánál készpénzben. Ebben az esetben egy emailt fogunk küldeni a konferencia előtt, melyben megkérünk majd egy második visszaigazolásra a megrendelt szolgáltatásokról."
$newText = "Amennyiben az utalást csak külföldi (nem magyar) bankszámláról tudnád intézni, úgy nem vagy köteles az előleg fizetésére tekintettel a magas határon kívüli utalási költségekre. A regisztráció költségét a konferencián személyesen, készpénzben tudod egyenlíteni. Ebben az esetben egy e-mailt fogunk küldeni a konferencia előtt, melyben megkérünk majd egy második visszaigazolásra a megrendelt szolgáltatásokról."
$d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
